$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "'4"
$ws.Range("D2").Value = 0.169
$ws.Range("E2").Value = 0.14
$ws.Range("G2").Value = 0.2690626517727052
$ws.Range("H2").Value = 0.2690626517727052
$ws.Range("I2").Value = 0.2611299983810911
$ws.Range("J2").Value = 0.2259874100803443
$ws.Range("K2").Value = 184
$ws.Range("L2").Value = 0.2978792293993848
$ws.Range("U2").Value = 663.1999999999999
$ws.Range("V2").Value = 0.5932552106628499
$ws.Range("W2").Value = 0.3298521823434459
$ws.Range("X2").Value = 0.05680993890095423
$ws.Range("Y2").Value = 0.2730422434424917
$ws.Range("Z2").Value = 1.104408336074836
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.0539982098105608
$ws.Range("AC2").Value = -0.04978686389103863
$ws.Range("AD2").Value = 411.003
$ws.Range("AF2").Value = 411.003
$ws.Range("AG2").Value = -252.1969999999999
$ws.Range("AH2").Value = 0.2688221554931869
$ws.Range("AI2").Value = 0.3299365900218592
$ws.Range("AJ2").Value = -0.2913204644086943
$ws.Range("AK2").Value = -0.432953993369991
$ws.Range("AL2").Value = 57.4
$ws.Range("AM2").Value = 57.4
$ws.Range("AN2").Value = 2.325993208828523
$ws.Range("AO2").Value = 2.810104529616725
$ws.Range("AP2").Value = -1.427260894170911
$ws.Range("AQ2").Value = 2.810104529616725
$ws.Range("D3").Value = 0.169
$ws.Range("E3").Value = 0.14
$ws.Range("G3").Value = 0.4934679334916864
$ws.Range("H3").Value = 0.4934679334916864
$ws.Range("I3").Value = 0.4789192399049881
$ws.Range("J3").Value = 0.4004835262403947
$ws.Range("K3").Value = 86.8
$ws.Range("L3").Value = 0.2577197149643705
$ws.Range("U3").Value = 500.9
$ws.Range("V3").Value = 1.066425377900788
$ws.Range("X3").Value = 0.06185898255388222
$ws.Range("AB3").Value = 0.05705564459784086
$ws.Range("AD3").Value = 210.6
$ws.Range("AF3").Value = 210.6
$ws.Range("AG3").Value = -290.3
$ws.Range("AH3").Value = 0.309569307658386
$ws.Range("AI3").Value = 0.3522328148519819
$ws.Range("AJ3").Value = -1.618171683389074
$ws.Range("AK3").Value = -2.992783505154637
$ws.Range("AL3").Value = 57.4
$ws.Range("AM3").Value = 57.4
$ws.Range("AN3").Value = 1.191850594227504
$ws.Range("AO3").Value = 2.810104529616725
$ws.Range("AP3").Value = -1.642897566496887
$ws.Range("AQ3").Value = 2.810104529616725
$ws.Range("B4").Value = 'Bank of Africa - Côte D''ivoire (BRVM:BOAC)'
$ws.Range("K4").Value = 25.7
$ws.Range("L4").Value = 0.4269102990033222
$ws.Range("U4").Value = 87.3
$ws.Range("V4").Value = 0.6367614879649891
$ws.Range("W4").Value = 0.3257287705956907
$ws.Range("X4").Value = 0.04863300287183572
$ws.Range("Y4").Value = 0.277095767723855
$ws.Range("Z4").Value = 4.892717815344603
$ws.Range("AB4").Value = 0.04863295275879653
$ws.Range("AC4").Value = -0.04863295275879653
$ws.Range("AD4").Value = 0.003
$ws.Range("AF4").Value = 0.003
$ws.Range("AG4").Value = -87.297
$ws.Range("AH4").Value = [double]"2.188135927003786e-05"
$ws.Range("AI4").Value = [double]"3.271430596599893e-05"
$ws.Range("AJ4").Value = -1.752846214083489
$ws.Range("AK4").Value = -19.82670906200315
$ws.Range("A5").Value = 'Ivory Coast'
$ws.Range("B5").Value = 'NSIA Banque Société Anonyme (BRVM:NSBC)'
$ws.Range("C5").Value = 'Bank (Money Center)'
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 19.5
$ws.Range("L5").Value = 0.1900584795321638
$ws.Range("M5").Value = -0
$ws.Range("N5").Value = -0
$ws.Range("O5").Value = -0
$ws.Range("P5").Value = -0
$ws.Range("Q5").Value = -0
$ws.Range("R5").Value = -0
$ws.Range("S5").Value = 0
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 0.1311365164761265
$ws.Range("X5").Value = 0.05176089524802623
$ws.Range("Y5").Value = 0.07937562122810023
$ws.Range("Z5").Value = 0.5354906054279749
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0.05094077502328073
$ws.Range("AC5").Value = -0.05094077502328073
$ws.Range("AD5").Value = 27.5
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 27.5
$ws.Range("AG5").Value = 27.5
$ws.Range("AH5").Value = 0.09588563458856346
$ws.Range("AI5").Value = 0.1355347461803844
$ws.Range("AJ5").Value = 0.09588563458856346
$ws.Range("AK5").Value = 0.1355347461803844
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = 0
$ws.Range("A6").Value = 'Ivory Coast'
$ws.Range("B6").Value = 'Société Ivoirienne de Banque S.A. (BRVM:SIBC)'
$ws.Range("C6").Value = 'Bank (Money Center)'
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 52
$ws.Range("L6").Value = 0.4403048264182896
$ws.Range("M6").Value = -0
$ws.Range("N6").Value = -0
$ws.Range("O6").Value = -0
$ws.Range("P6").Value = -0
$ws.Range("Q6").Value = -0
$ws.Range("R6").Value = -0
$ws.Range("S6").Value = 0
$ws.Range("U6").Value = 75
$ws.Range("V6").Value = 0.2978554408260524
$ws.Range("W6").Value = 0.333975594091201
$ws.Range("X6").Value = 0.068888209011538
$ws.Range("Y6").Value = 0.265087385079663
$ws.Range("Z6").Value = 0.3323016319639843
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 0.0597097218708334
$ws.Range("AC6").Value = -0.0597097218708334
$ws.Range("AD6").Value = 172.9
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 172.9
$ws.Range("AG6").Value = 97.90000000000001
$ws.Range("AH6").Value = 0.4071109018130445
$ws.Range("AI6").Value = 0.489524348810872
$ws.Range("AJ6").Value = 0.2799542464969974
$ws.Range("AK6").Value = 0.3519051042415528
$ws.Range("AL6").Value = 0
$ws.Range("AM6").Value = 0
